$wb = $excel.ActiveWorkbook

# Sheet "展览" (sheetId 1)
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 7775
$ws1.Range("F3").Value = 73
$ws1.Range("F5").Value = 54
$ws1.Range("F6").Value = 510
$ws1.Range("F7").Value = 1160
$ws1.Range("F8").Value = 207
$ws1.Range("F10").Value = 167
$ws1.Range("F11").Value = 40

# Sheet "全部类型" (sheetId 4)
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 7775
$ws4.Range("F3").Value = 73
$ws4.Range("F5").Value = 54
$ws4.Range("F6").Value = 510
$ws4.Range("F7").Value = 1160
$ws4.Range("F8").Value = 207
$ws4.Range("F11").Value = 167
$ws4.Range("F12").Value = 40
